# Update recomputed TPM-derived NATMI metrics (ligand/receptor expression,
# specificity, and edge-weight columns) for the Thbs1-Itga3 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target range G2:T26 -> 25 rows x 14 cols (columns G..T). The K/L rank
# columns are included in the range but keep their original values.
$arr = New-Object 'object[,]' 25,14

# Row 2
$arr[0,0] = 41.428665
$arr[0,1] = 124.285995
$arr[0,2] = 0.06969137269740189
$arr[0,3] = 0.06969137269740189
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 15.47987166666667
$arr[0,7] = 46.439615
$arr[0,8] = 0.7960757698994193
$arr[0,9] = 0.7960757698994194
$arr[0,10] = 641.310417521325
$arr[0,11] = 5771.793757691925
$arr[0,12] = 0.05547961317543158
$arr[0,13] = 0.05547961317543158
# Row 3
$arr[1,0] = 41.428665
$arr[1,1] = 124.285995
$arr[1,2] = 0.06969137269740189
$arr[1,3] = 0.06969137269740189
$arr[1,4] = 3
$arr[1,5] = 1
$arr[1,6] = 0.9918089999999999
$arr[1,7] = 2.975427
$arr[1,8] = 0.05100527512565552
$arr[1,9] = 0.05100527512565553
$arr[1,10] = 41.089322804985
$arr[1,11] = 369.8039052448649
$arr[1,12] = 0.003554627638315581
$arr[1,13] = 0.003554627638315581
# Row 4
$arr[2,0] = 41.428665
$arr[2,1] = 124.285995
$arr[2,2] = 0.06969137269740189
$arr[2,3] = 0.06969137269740189
$arr[2,4] = 2
$arr[2,5] = 0.6666666666666666
$arr[2,6] = 0.171678
$arr[2,7] = 0.515034
$arr[2,8] = 0.008828800326496623
$arr[2,9] = 0.008828800326496624
$arr[2,10] = 7.112390349869999
$arr[2,11] = 64.01151314882999
$arr[2,12] = 0.0006152912140248196
$arr[2,13] = 0.0006152912140248197
# Row 5
$arr[3,0] = 41.428665
$arr[3,1] = 124.285995
$arr[3,2] = 0.06969137269740189
$arr[3,3] = 0.06969137269740189
$arr[3,4] = 3
$arr[3,5] = 1
$arr[3,6] = 2.688466333333334
$arr[3,7] = 8.065399000000001
$arr[3,8] = 0.1382584398787761
$arr[3,9] = 0.1382584398787761
$arr[3,10] = 111.379571087445
$arr[3,11] = 1002.416139787005
$arr[3,12] = 0.009635420462153115
$arr[3,13] = 0.009635420462153115
# Row 6
$arr[4,0] = 41.428665
$arr[4,1] = 124.285995
$arr[4,2] = 0.06969137269740189
$arr[4,3] = 0.06969137269740189
$arr[4,4] = 1
$arr[4,5] = 0.3333333333333333
$arr[4,6] = 0.113399
$arr[4,7] = 0.340197
$arr[4,8] = 0.005831714769652435
$arr[4,9] = 0.005831714769652436
$arr[4,10] = 4.697969182335
$arr[4,11] = 42.281722641015
$arr[4,12] = 0.000406420207476791
$arr[4,13] = 0.0004064202074767911
# Row 7
$arr[5,0] = 170.232249
$arr[5,1] = 510.696747
$arr[5,2] = 0.2863649869040173
$arr[5,3] = 0.2863649869040173
$arr[5,4] = 3
$arr[5,5] = 1
$arr[5,6] = 15.47987166666667
$arr[5,7] = 46.439615
$arr[5,8] = 0.7960757698994193
$arr[5,9] = 0.7960757698994194
$arr[5,10] = 2635.173368048045
$arr[5,11] = 23716.5603124324
$arr[5,12] = 0.2279682274218527
$arr[5,13] = 0.2279682274218528
# Row 8
$arr[6,0] = 170.232249
$arr[6,1] = 510.696747
$arr[6,2] = 0.2863649869040173
$arr[6,3] = 0.2863649869040173
$arr[6,4] = 3
$arr[6,5] = 1
$arr[6,6] = 0.9918089999999999
$arr[6,7] = 2.975427
$arr[6,8] = 0.05100527512565552
$arr[6,9] = 0.05100527512565553
$arr[6,10] = 168.837876648441
$arr[6,11] = 1519.540889835969
$arr[6,12] = 0.01460612494339414
$arr[6,13] = 0.01460612494339415
# Row 9
$arr[7,0] = 170.232249
$arr[7,1] = 510.696747
$arr[7,2] = 0.2863649869040173
$arr[7,3] = 0.2863649869040173
$arr[7,4] = 2
$arr[7,5] = 0.6666666666666666
$arr[7,6] = 0.171678
$arr[7,7] = 0.515034
$arr[7,8] = 0.008828800326496623
$arr[7,9] = 0.008828800326496624
$arr[7,10] = 29.225132043822
$arr[7,11] = 263.026188394398
$arr[7,12] = 0.002528259289875389
$arr[7,13] = 0.00252825928987539
# Row 10
$arr[8,0] = 170.232249
$arr[8,1] = 510.696747
$arr[8,2] = 0.2863649869040173
$arr[8,3] = 0.2863649869040173
$arr[8,4] = 3
$arr[8,5] = 1
$arr[8,6] = 2.688466333333334
$arr[8,7] = 8.065399000000001
$arr[8,8] = 0.1382584398787761
$arr[8,9] = 0.1382584398787761
$arr[8,10] = 457.6636702841171
$arr[8,11] = 4118.973032557054
$arr[8,12] = 0.03959237632525557
$arr[8,13] = 0.03959237632525558
# Row 11
$arr[9,0] = 170.232249
$arr[9,1] = 510.696747
$arr[9,2] = 0.2863649869040173
$arr[9,3] = 0.2863649869040173
$arr[9,4] = 1
$arr[9,5] = 0.3333333333333333
$arr[9,6] = 0.113399
$arr[9,7] = 0.340197
$arr[9,8] = 0.005831714769652435
$arr[9,9] = 0.005831714769652436
$arr[9,10] = 19.304166804351
$arr[9,11] = 173.737501239159
$arr[9,12] = 0.001669998923639484
$arr[9,13] = 0.001669998923639484
# Row 12
$arr[10,0] = 244.5761666666666
$arr[10,1] = 733.7284999999999
$arr[10,2] = 0.4114264551867299
$arr[10,3] = 0.41142645518673
$arr[10,4] = 3
$arr[10,5] = 1
$arr[10,6] = 15.47987166666667
$arr[10,7] = 46.439615
$arr[10,8] = 0.7960757698994193
$arr[10,9] = 0.7960757698994194
$arr[10,10] = 3786.007672725278
$arr[10,11] = 34074.0690545275
$arr[10,12] = 0.327526632069765
$arr[10,13] = 0.3275266320697651
# Row 13
$arr[11,0] = 244.5761666666666
$arr[11,1] = 733.7284999999999
$arr[11,2] = 0.4114264551867299
$arr[11,3] = 0.41142645518673
$arr[11,4] = 3
$arr[11,5] = 1
$arr[11,6] = 0.9918089999999999
$arr[11,7] = 2.975427
$arr[11,8] = 0.05100527512565552
$arr[11,9] = 0.05100527512565553
$arr[11,10] = 242.5728432854999
$arr[11,11] = 2183.1555895695
$arr[11,12] = 0.02098491954077234
$arr[11,13] = 0.02098491954077235
# Row 14
$arr[12,0] = 244.5761666666666
$arr[12,1] = 733.7284999999999
$arr[12,2] = 0.4114264551867299
$arr[12,3] = 0.41142645518673
$arr[12,4] = 2
$arr[12,5] = 0.6666666666666666
$arr[12,6] = 0.171678
$arr[12,7] = 0.515034
$arr[12,8] = 0.008828800326496623
$arr[12,9] = 0.008828800326496624
$arr[12,10] = 41.98834714099999
$arr[12,11] = 377.895124269
$arr[12,12] = 0.003632402021881949
$arr[12,13] = 0.003632402021881951
# Row 15
$arr[13,0] = 244.5761666666666
$arr[13,1] = 733.7284999999999
$arr[13,2] = 0.4114264551867299
$arr[13,3] = 0.41142645518673
$arr[13,4] = 3
$arr[13,5] = 1
$arr[13,6] = 2.688466333333334
$arr[13,7] = 8.065399000000001
$arr[13,8] = 0.1382584398787761
$arr[13,9] = 0.1382584398787761
$arr[13,10] = 657.5347900190555
$arr[13,11] = 5917.813110171501
$arr[13,12] = 0.05688317981897245
$arr[13,13] = 0.05688317981897247
# Row 16
$arr[14,0] = 244.5761666666666
$arr[14,1] = 733.7284999999999
$arr[14,2] = 0.4114264551867299
$arr[14,3] = 0.41142645518673
$arr[14,4] = 1
$arr[14,5] = 0.3333333333333333
$arr[14,6] = 0.113399
$arr[14,7] = 0.340197
$arr[14,8] = 0.005831714769652435
$arr[14,9] = 0.005831714769652436
$arr[14,10] = 27.73469272383333
$arr[14,11] = 249.6122345145
$arr[14,12] = 0.002399321735338199
$arr[14,13] = 0.002399321735338199
# Row 17
$arr[15,0] = 24.173247
$arr[15,1] = 72.51974100000001
$arr[15,2] = 0.04066427836821081
$arr[15,3] = 0.04066427836821081
$arr[15,4] = 3
$arr[15,5] = 1
$arr[15,6] = 15.47987166666667
$arr[15,7] = 46.439615
$arr[15,8] = 0.7960757698994193
$arr[15,9] = 0.7960757698994194
$arr[15,10] = 374.1987613266351
$arr[15,11] = 3367.788851939716
$arr[15,12] = 0.03237184670937773
$arr[15,13] = 0.03237184670937773
# Row 18
$arr[16,0] = 24.173247
$arr[16,1] = 72.51974100000001
$arr[16,2] = 0.04066427836821081
$arr[16,3] = 0.04066427836821081
$arr[16,4] = 3
$arr[16,5] = 1
$arr[16,6] = 0.9918089999999999
$arr[16,7] = 2.975427
$arr[16,8] = 0.05100527512565552
$arr[16,9] = 0.05100527512565553
$arr[16,10] = 23.975243933823
$arr[16,11] = 215.777195404407
$arr[16,12] = 0.002074092705956835
$arr[16,13] = 0.002074092705956835
# Row 19
$arr[17,0] = 24.173247
$arr[17,1] = 72.51974100000001
$arr[17,2] = 0.04066427836821081
$arr[17,3] = 0.04066427836821081
$arr[17,4] = 2
$arr[17,5] = 0.6666666666666666
$arr[17,6] = 0.171678
$arr[17,7] = 0.515034
$arr[17,8] = 0.008828800326496623
$arr[17,9] = 0.008828800326496624
$arr[17,10] = 4.150014698466
$arr[17,11] = 37.350132286194
$arr[17,12] = 0.0003590167941340092
$arr[17,13] = 0.0003590167941340092
# Row 20
$arr[18,0] = 24.173247
$arr[18,1] = 72.51974100000001
$arr[18,2] = 0.04066427836821081
$arr[18,3] = 0.04066427836821081
$arr[18,4] = 3
$arr[18,5] = 1
$arr[18,6] = 2.688466333333334
$arr[18,7] = 8.065399000000001
$arr[18,8] = 0.1382584398787761
$arr[18,9] = 0.1382584398787761
$arr[18,10] = 64.98896072685102
$arr[18,11] = 584.9006465416592
$arr[18,12] = 0.005622179685985089
$arr[18,13] = 0.005622179685985089
# Row 21
$arr[19,0] = 24.173247
$arr[19,1] = 72.51974100000001
$arr[19,2] = 0.04066427836821081
$arr[19,3] = 0.04066427836821081
$arr[19,4] = 1
$arr[19,5] = 0.3333333333333333
$arr[19,6] = 0.113399
$arr[19,7] = 0.340197
$arr[19,8] = 0.005831714769652435
$arr[19,9] = 0.005831714769652436
$arr[19,10] = 2.741222036553001
$arr[19,11] = 24.67099832897701
$arr[19,12] = 0.000237142472757153
$arr[19,13] = 0.000237142472757153
# Row 22
$arr[20,0] = 114.0486906666667
$arr[20,1] = 342.146072
$arr[20,2] = 0.19185290684364
$arr[20,3] = 0.19185290684364
$arr[20,4] = 3
$arr[20,5] = 1
$arr[20,6] = 15.47987166666667
$arr[20,7] = 46.439615
$arr[20,8] = 0.7960757698994193
$arr[20,9] = 0.7960757698994194
$arr[20,10] = 1765.459095271365
$arr[20,11] = 15889.13185744228
$arr[20,12] = 0.1527294505229923
$arr[20,13] = 0.1527294505229923
# Row 23
$arr[21,0] = 114.0486906666667
$arr[21,1] = 342.146072
$arr[21,2] = 0.19185290684364
$arr[21,3] = 0.19185290684364
$arr[21,4] = 3
$arr[21,5] = 1
$arr[21,6] = 0.9918089999999999
$arr[21,7] = 2.975427
$arr[21,8] = 0.05100527512565552
$arr[21,9] = 0.05100527512565553
$arr[21,10] = 113.114517841416
$arr[21,11] = 1018.030660572744
$arr[21,12] = 0.009785510297216617
$arr[21,13] = 0.009785510297216618
# Row 24
$arr[22,0] = 114.0486906666667
$arr[22,1] = 342.146072
$arr[22,2] = 0.19185290684364
$arr[22,3] = 0.19185290684364
$arr[22,4] = 2
$arr[22,5] = 0.6666666666666666
$arr[22,6] = 0.171678
$arr[22,7] = 0.515034
$arr[22,8] = 0.008828800326496623
$arr[22,9] = 0.008828800326496624
$arr[22,10] = 19.579651116272
$arr[22,11] = 176.216860046448
$arr[22,12] = 0.001693831006580455
$arr[22,13] = 0.001693831006580455
# Row 25
$arr[23,0] = 114.0486906666667
$arr[23,1] = 342.146072
$arr[23,2] = 0.19185290684364
$arr[23,3] = 0.19185290684364
$arr[23,4] = 3
$arr[23,5] = 1
$arr[23,6] = 2.688466333333334
$arr[23,7] = 8.065399000000001
$arr[23,8] = 0.1382584398787761
$arr[23,9] = 0.1382584398787761
$arr[23,10] = 306.616065218081
$arr[23,11] = 2759.544586962728
$arr[23,12] = 0.02652528358640982
$arr[23,13] = 0.02652528358640982
# Row 26
$arr[24,0] = 114.0486906666667
$arr[24,1] = 342.146072
$arr[24,2] = 0.19185290684364
$arr[24,3] = 0.19185290684364
$arr[24,4] = 1
$arr[24,5] = 0.3333333333333333
$arr[24,6] = 0.113399
$arr[24,7] = 0.340197
$arr[24,8] = 0.005831714769652435
$arr[24,9] = 0.005831714769652436
$arr[24,10] = 12.93300747290933
$arr[24,11] = 116.397067256184
$arr[24,12] = 0.001118831430440808
$arr[24,13] = 0.001118831430440808

$ws.Range("G2:T26").Value = $arr
